$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in G1 (match header row style: bold, centered - same as A1:F1)
$ws.Range("G1").Value = "Completion Date"
$ws.Range("G1").Font.Bold = $true
$ws.Range("G1").HorizontalAlignment = -4108
$ws.Range("G1").VerticalAlignment = -4108

# Set column G width to match the header's best fit width (~15.21875, closest achievable)
$ws.Columns("G").ColumnWidth = 14.3

# Add completion dates for rows 2 and 7 (value 43379 -> 2018-10-06)
$ws.Range("G2").Value = 43379
$ws.Range("G2").NumberFormat = "d-mmm-yy"
$ws.Range("G2").Font.Color = $ws.Range("A2").Font.Color

$ws.Range("G7").Value = 43379
$ws.Range("G7").NumberFormat = "d-mmm-yy"
$ws.Range("G7").Font.Color = $ws.Range("A7").Font.Color

# Update selection
$ws.Range("G3").Select()
